$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.250.63'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '2.550.31'
$ws.Range('E3').Value = '  +3.28%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.11%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('D9').Value = '2.547.26'
$ws.Range('E9').Value = '  +3.22%  '
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('E14').Value = '  +3.98%  '
$ws.Range('D15').Value = '3.005.52'
$ws.Range('E15').Value = '  +3.34%  '
$ws.Range('D16').Value = '63.186.96'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('D18').Value = '2.551.63'
$ws.Range('E18').Value = '  +3.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +8.79%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('E27').Value = '  +11.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.90%  '
$ws.Range('D31').Value = '0.0₃0824'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '178.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '412.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.40%  '
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('E40').Value = '  +4.32%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.56'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.29%  '
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0240'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.08%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0522'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.92%  '
$ws.Range('E51').Value = '  +2.41%  '
